# Self-Assessment workbook update:
# Student 202007021 (row 12) graded 5/5/5/5 across the four evaluators
# (columns D:G). The AVERAGE() formulas in column H (per-student) and
# row 14 (per-evaluator) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Group and Self Assessment")

# Fill in the grades for the row-12 student (TeamID 202007021).
$ws.Range("D12:G12").Value = 5

# Leave the selection where the author last clicked.
$ws.Range("J12").Select() | Out-Null
